$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2261.2222
$ws.Range("I40").Value = 1817
$ws.Range("J40").Value = 2483.3333
$ws.Range("K40").Value = 1817
$ws.Range("L40").Value = 2483.3333
$ws.Range("M40").Value = -1642
$ws.Range("N40").Value = -2833.3333

$ws.Range("H106").Value = 1616.6666
$ws.Range("I106").Value = 1540
$ws.Range("J106").Value = 2000
$ws.Range("K106").Value = 1540
$ws.Range("L106").Value = 2000
$ws.Range("M106").Value = -909
$ws.Range("N106").Value = -3262

$ws.Range("H115").Value = 2428.5715
$ws.Range("J115").Value = 2428.5715
$ws.Range("L115").Value = 7285.7145
$ws.Range("N115").Value = -10419.7145

$ws.Range("H129").Value = 3522353.5
$ws.Range("J129").Value = 1238.1406
$ws.Range("L129").Value = 3714.4218
$ws.Range("N129").Value = -13714.4218

$ws.Range("H135").Value = 1363.8235
$ws.Range("I135").Value = 1234.2142
$ws.Range("J135").Value = 1968.6666
$ws.Range("K135").Value = 11107.9278
$ws.Range("L135").Value = 17717.9994
$ws.Range("M135").Value = -8572.927799999999
$ws.Range("N135").Value = -22787.9994

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 16668600
$ws.Range("I132").Value = 20409458
$ws.Range("J132").Value = 4773.8184
$ws.Range("K132").Value = 61228374
$ws.Range("L132").Value = 14321.4552
$ws.Range("M132").Value = -61225844
$ws.Range("N132").Value = -19381.4552

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H24").Value = 24520.555
$ws.Range("I24").Value = 1771.8334
$ws.Range("J24").Value = 70018
$ws.Range("K24").Value = 1771.8334
$ws.Range("L24").Value = 70018
$ws.Range("M24").Value = -1536.8334
$ws.Range("N24").Value = -70488

$ws.Range("H51").Value = 30774.117
$ws.Range("J51").Value = 30774.117
$ws.Range("L51").Value = 30774.117
$ws.Range("N51").Value = -31756.117

$ws.Range("H134").Value = 2420.2354
$ws.Range("I134").Value = 1773.3572
$ws.Range("J134").Value = 5439
$ws.Range("K134").Value = 5320.071599999999
$ws.Range("L134").Value = 16317
$ws.Range("M134").Value = -2785.071599999999
$ws.Range("N134").Value = -21387

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6179.4116
$ws.Range("I3").Value = 3905.6
$ws.Range("J3").Value = 7126.8335
$ws.Range("K3").Value = 11716.8
$ws.Range("L3").Value = 21380.5005
$ws.Range("M3").Value = -11604.8
$ws.Range("N3").Value = -21604.5005

$ws.Range("H37").Value = 70000
$ws.Range("J37").Value = 70000
$ws.Range("L37").Value = 210000
$ws.Range("N37").Value = -210224

$ws.Range("H74").Value = 14230
$ws.Range("I74").Value = 1000
$ws.Range("J74").Value = 15883.75
$ws.Range("K74").Value = 3000
$ws.Range("L74").Value = 47651.25
$ws.Range("M74").Value = -1939
$ws.Range("N74").Value = -49773.25

$ws.Range("H75").Value = 3000
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 3000
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 9000
$ws.Range("M75").Value = ""
$ws.Range("N75").Value = -10996

$ws.Range("H76").Value = 3493.2222
$ws.Range("I76").Value = 2679.6667
$ws.Range("J76").Value = 3900
$ws.Range("K76").Value = 8039.000100000001
$ws.Range("L76").Value = 11700
$ws.Range("M76").Value = -7656.000100000001
$ws.Range("N76").Value = -12466

$ws.Range("H77").Value = 14230
$ws.Range("I77").Value = 1000
$ws.Range("J77").Value = 15883.75
$ws.Range("K77").Value = 9000
$ws.Range("L77").Value = 142953.75
$ws.Range("M77").Value = -3696
$ws.Range("N77").Value = -153561.75

$ws.Range("H78").Value = 3000
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 3000
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 27000
$ws.Range("M78").Value = ""
$ws.Range("N78").Value = -36984

$ws.Range("H79").Value = 3493.2222
$ws.Range("I79").Value = 2679.6667
$ws.Range("J79").Value = 3900
$ws.Range("K79").Value = 8039.000100000001
$ws.Range("L79").Value = 11700
$ws.Range("M79").Value = -6713.000100000001
$ws.Range("N79").Value = -14352

$ws.Range("H113").Value = 7693355
$ws.Range("I113").Value = 100000000
$ws.Range("J113").Value = 1134.5
$ws.Range("K113").Value = 300000000
$ws.Range("L113").Value = 3403.5
$ws.Range("M113").Value = -299997830
$ws.Range("N113").Value = -7743.5

$ws.Range("H133").Value = 7171.6665
$ws.Range("I133").Value = 8757.5
$ws.Range("K133").Value = 26272.5
$ws.Range("M133").Value = -21212.5

$ws.Range("H134").Value = 2713.611
$ws.Range("I134").Value = 1654
$ws.Range("J134").Value = 3387.9092
$ws.Range("K134").Value = 4962
$ws.Range("L134").Value = 10163.7276
$ws.Range("M134").Value = 108
$ws.Range("N134").Value = -20303.7276

$ws.Range("H136").Value = 2467.0625
$ws.Range("I136").Value = 1789.091
$ws.Range("J136").Value = 3958.6
$ws.Range("K136").Value = 5367.272999999999
$ws.Range("L136").Value = 11875.8
$ws.Range("M136").Value = -267.2729999999992
$ws.Range("N136").Value = -22075.8

$ws.Range("H137").Value = 2598.7856
$ws.Range("I137").Value = 1983.3334
$ws.Range("J137").Value = 3706.6
$ws.Range("K137").Value = 5950.0002
$ws.Range("L137").Value = 11119.8
$ws.Range("M137").Value = -850.0002000000004
$ws.Range("N137").Value = -21319.8

$ws.Range("H139").Value = 8163
$ws.Range("I139").Value = 2230.5
$ws.Range("J139").Value = 22005.5
$ws.Range("K139").Value = 6691.5
$ws.Range("L139").Value = 66016.5
$ws.Range("M139").Value = -1551.5
$ws.Range("N139").Value = -76296.5

$ws.Range("H140").Value = 3821.35
$ws.Range("I140").Value = 995.4
$ws.Range("J140").Value = 4763.3335
$ws.Range("K140").Value = 2986.2
$ws.Range("L140").Value = 14290.0005
$ws.Range("M140").Value = 2193.8
$ws.Range("N140").Value = -24650.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H119").Value = 35860
$ws.Range("J119").Value = 35860
$ws.Range("L119").Value = 35860
$ws.Range("N119").Value = -45536

$ws.Range("H132").Value = 2956.8628
$ws.Range("I132").Value = 2568.4062
$ws.Range("J132").Value = 3611.1052
$ws.Range("K132").Value = 7705.2186
$ws.Range("L132").Value = 10833.3156
$ws.Range("M132").Value = -5175.2186
$ws.Range("N132").Value = -15893.3156

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1721.5555
$ws.Range("I7").Value = 1618.8
$ws.Range("J7").Value = 1850
$ws.Range("K7").Value = 1618.8
$ws.Range("L7").Value = 1850
$ws.Range("M7").Value = -1506.8
$ws.Range("N7").Value = -2074

$ws.Range("H46").Value = 2220.9167
$ws.Range("J46").Value = 3380
$ws.Range("L46").Value = 3380
$ws.Range("N46").Value = -3756

$ws.Range("H122").Value = 2931.5
$ws.Range("I122").Value = 2223.4443
$ws.Range("J122").Value = 4524.625
$ws.Range("K122").Value = 6670.3329
$ws.Range("L122").Value = 13573.875
$ws.Range("M122").Value = -4220.3329
$ws.Range("N122").Value = -18473.875

$ws.Range("H126").Value = 1721.5555
$ws.Range("I126").Value = 1618.8
$ws.Range("J126").Value = 1850
$ws.Range("K126").Value = 4856.4
$ws.Range("L126").Value = 5550
$ws.Range("M126").Value = -2386.4
$ws.Range("N126").Value = -10490

$ws.Range("H132").Value = 3055.6453
$ws.Range("I132").Value = 1707.1765
$ws.Range("J132").Value = 4693.0713
$ws.Range("K132").Value = 5121.529500000001
$ws.Range("L132").Value = 14079.2139
$ws.Range("M132").Value = -2591.529500000001
$ws.Range("N132").Value = -19139.2139

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 51012.668
$ws.Range("J31").Value = 51012.668
$ws.Range("L31").Value = 51012.668
$ws.Range("N31").Value = -51708.668

$ws.Range("H122").Value = 669383.8
$ws.Range("I122").Value = 835396.5
$ws.Range("K122").Value = 2506189.5
$ws.Range("M122").Value = -2503739.5

$ws.Range("H127").Value = 38333.332
$ws.Range("J127").Value = 38333.332
$ws.Range("L127").Value = 38333.332
$ws.Range("N127").Value = -48253.332
